# Generate Report for Handoff
# Update the "Latest Handoff Datetime" for the 74a9efa0-875e-445a-a488-ab5e820d6ca8
# file row (row 5) on both the zh-cn and de-de status sheets to reflect a new
# handoff that just occurred.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-09 12:33:56"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-09 12:34:02"
